$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: fix values that the fuzzer had shifted/corrupted ---
$ws.Range("G12").Value2 = 1240524717.0500007
$ws.Range("G13").Value2 = 319819483.18000001
$ws.Range("G14").Value2 = 34063116.800000042
$ws.Range("G15").Value2 = 40000000
$ws.Range("G16").Value2 = -60834434.380000003

# G18 becomes a real SUM formula (like the other columns) instead of a hard-coded value
$ws.Range("G18").Formula = "=SUM(G12:G17)"

$ws.Range("G19").Value2 = -379300000.00000012

# G21 becomes a real SUM formula (like the other columns) instead of a hard-coded value
$ws.Range("G21").Formula = "=SUM(G18:G20)"

$ws.Range("G22").Value2 = -20015625

# G26 was mistakenly a shared-string blank; make it the correct numeric value
$ws.Range("G26").Value2 = 1029174575.116062

# --- Column I: correct minor rounding / data-entry glitches ---
$ws.Range("I12").Value2 = 4188377156
$ws.Range("I13").Value2 = 1012006300
$ws.Range("I14").Value2 = -44319159.289999999
$ws.Range("I16").Value2 = -162861893.59999999
$ws.Range("I19").Value2 = -1160500000
$ws.Range("I26").Value2 = 10110658959

Write-Output "done"
